$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (row, date-serial, B, C, D) appended after existing row 357
$rows = @(
    @(358, 44432, 0, 4, 60.91989034419738),
    @(359, 44433, 0, 3, 45.68991775814803),
    @(360, 44434, 1, 4, 60.91989034419738),
    @(361, 44435, 2, 6, 91.37983551629607),
    @(362, 44436, 2, 8, 121.8397806883948),
    @(363, 44437, 2, 7, 106.6098081023454),
    @(364, 44438, 7, 14, 213.2196162046908),
    @(365, 44439, 1, 15, 228.4495887907402),
    @(366, 44440, 0, 15, 228.4495887907402)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]

    # Column A: date serial value, carrying over the same formatted style as the
    # preceding row's date cell (centered/top, bordered, custom date format).
    $ws.Cells.Item($rowIndex - 1, 1).Copy($ws.Cells.Item($rowIndex, 1))
    $ws.Cells.Item($rowIndex, 1).Value = $r[1]

    # Columns B, C, D: plain numbers, unstyled like the rest of the table
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
}
